$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 10; this shifts existing rows 10..90 down to 11..91
$ws.Rows.Item(10).Insert()

# Copy formatting (style) of the date cell from the row that used to be row 10
# (now row 11) onto the new row's date cell so the date style (s="2") matches.
$ws.Range("D11").Copy() | Out-Null
$ws.Range("D10").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# Fill in the new row's data
$ws.Cells.Item(10, 1).Value = 11
$ws.Cells.Item(10, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(10, 3).Value = "Bíobío"
$ws.Cells.Item(10, 4).Value = 45168
$ws.Cells.Item(10, 5).Value = 8
$ws.Cells.Item(10, 6).Value = 100112031
$ws.Cells.Item(10, 7).Value = "Poroto verde"
$ws.Cells.Item(10, 8).Value = "Sin especificar"
$ws.Cells.Item(10, 9).Value = "Primera"
$ws.Cells.Item(10, 10).Value = 40
$ws.Cells.Item(10, 11).Value = 24000
$ws.Cells.Item(10, 12).Value = 24000
$ws.Cells.Item(10, 13).Value = 24000
$ws.Cells.Item(10, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(10, 15).Value = "Perú"
$ws.Cells.Item(10, 16).Value = 960
$ws.Cells.Item(10, 17).Value = 25
$ws.Cells.Item(10, 18).Value = "Hortaliza"
